$wb = $excel.ActiveWorkbook

# --- Update shared text labels (column headers) on both sheets ---
$wsMeans = $wb.Worksheets.Item("Means")
$wsSD = $wb.Worksheets.Item("Standard Deviations")

$wsMeans.Range("B1").Value = "National Average"
$wsMeans.Range("C1").Value = "State Average"

$wsSD.Range("B1").Value = "National Average SD"
$wsSD.Range("C1").Value = "State Average SD"

# --- Update values in "Means" sheet (sheet1) column B ---
$wsMeans.Range("B2").Value = 72
$wsMeans.Range("B3").Value = 13
$wsMeans.Range("B4").Value = 15
$wsMeans.Range("B5").Value = 18
$wsMeans.Range("B6").Value = 71
$wsMeans.Range("B7").Value = 7.3
$wsMeans.Range("B8").Value = 5.8
$wsMeans.Range("B9").Value = 29
$wsMeans.Range("B10").Value = 0.37

# --- Update values in "Standard Deviations" sheet (sheet2) column B ---
$wsSD.Range("B2").Value = 27
$wsSD.Range("B3").Value = 23
$wsSD.Range("B4").Value = 16
$wsSD.Range("B5").Value = 22
$wsSD.Range("B6").Value = 37
$wsSD.Range("B7").Value = 8.7
$wsSD.Range("B8").Value = 7.8
$wsSD.Range("B9").Value = 10
